# 11.24 batch-temp.md
# - CONFIG: clear the stray empty value cell at B11 (defaults.state_fallback row)
#   so the cell is dropped entirely, matching the other rows that have no
#   default value populated yet.
# - INPUT_MASTER: rename the header row to the BD_-prefixed column names.
# NOTE: config keys (CONFIG!A:A) and blacklist names (BLACKLIST_NAMES!A:A)
# are intentionally left untouched.

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("CONFIG")
$wsConfig.Range("B11").ClearContents()

$wsInput = $wb.Worksheets.Item("INPUT_MASTER")
$headerMap = @{
    "A1" = "BD_RECORD_ID"
    "B1" = "BD_SOURCE_TYPE"
    "C1" = "BD_ENTITY_NAME"
    "D1" = "BD_SOURCE_ENTITY_ID"
    "E1" = "BD_TITLE_ROLE"
    "F1" = "BD_TARGET_FIRST_NAME"
    "G1" = "BD_TARGET_LAST_NAME"
    "H1" = "BD_OWNER_NAME_FULL"
    "I1" = "BD_ADDRESS"
    "J1" = "BD_ADDRESS_2"
    "K1" = "BD_CITY"
    "L1" = "BD_STATE"
    "M1" = "BD_ZIP"
    "N1" = "BD_COUNTY"
    "O1" = "BD_APN"
    "P1" = "BD_MAILING_LINE1"
    "Q1" = "BD_MAILING_CITY"
    "R1" = "BD_MAILING_STATE"
    "S1" = "BD_MAILING_ZIP"
    "T1" = "BD_NOTES"
}

foreach ($addr in $headerMap.Keys) {
    $wsInput.Range($addr).Value = $headerMap[$addr]
}
